# Apply cryptos list update (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.446.73'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '3.125.89'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.51'
$ws.Range('E5').Value = '  +8.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '634.94'
$ws.Range('E6').Value = '  +1.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.08'
$ws.Range('E7').Value = '  +11.12%  '
$ws.Range('E8').Value = '  -6.50%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').Value = '3.124.48'
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.727'
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('E12').Value = '  +4.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.43'
$ws.Range('E13').Value = '  +5.23%  '
$ws.Range('B14').Value = 'Toncoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.63'
$ws.Range('E14').Value = '  +4.53%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000243'
$ws.Range('E15').Value = '  -4.69%  '
$ws.Range('D16').Value = '90.483.68'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').Value = '3.702.92'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '3.116.90'
$ws.Range('E18').Value = '  +0.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.65'
$ws.Range('E19').Value = '  -2.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.40'
$ws.Range('E20').Value = '  +2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000211'
$ws.Range('E21').Value = '  -4.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '450.23'
$ws.Range('E22').Value = '  +2.81%  '
$ws.Range('E23').Value = '  +10.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.06'
$ws.Range('E24').Value = '  +3.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.96'
$ws.Range('E25').Value = '  -3.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '90.73'
$ws.Range('E26').Value = '  +4.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.41'
$ws.Range('E27').Value = '  +1.42%  '
$ws.Range('D28').Value = '3.299.98'
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.99'
$ws.Range('E30').Value = '  +9.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.160'
$ws.Range('E31').Value = '  -4.92%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.204'
$ws.Range('E32').Value = '  +33.91%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.30'
$ws.Range('E33').Value = '  +15.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.84'
$ws.Range('E34').Value = '  +2.12%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.150'
$ws.Range('E35').Value = '  +6.38%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '513.09'
$ws.Range('E36').Value = '  -2.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.11'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.93'
$ws.Range('E38').Value = '  +4.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.31'
$ws.Range('E39').Value = '  +2.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.802'
$ws.Range('E40').Value = '  -19.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.424'
$ws.Range('E41').Value = '  +12.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0869'
$ws.Range('E42').Value = '  +3.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.18'
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('E45').Value = '  +43.45%  '
$ws.Range('E46').Value = '  +2.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.702'
$ws.Range('E47').Value = '  +13.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '150.21'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.61'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.26'
$ws.Range('E50').Value = '  +2.86%  '
$ws.Range('E51').Value = '  +3.93%  '
